$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 6).Value = 895
$ws.Cells.Item(2, 7).Value = 70
$ws.Cells.Item(3, 6).Value = 13957
$ws.Cells.Item(3, 7).Value = 88
$ws.Cells.Item(4, 6).Value = 13790
$ws.Cells.Item(4, 7).Value = 90
$ws.Cells.Item(5, 6).Value = 1065
$ws.Cells.Item(5, 7).Value = 70
$ws.Cells.Item(7, 7).Value = "不可售"
$ws.Cells.Item(9, 7).Value = "不可售"
$ws.Cells.Item(12, 6).Value = 784
$ws.Cells.Item(14, 6).Value = 142
$ws.Cells.Item(15, 6).Value = 105
$ws.Cells.Item(16, 6).Value = 89
$ws.Cells.Item(17, 6).Value = 154
$ws.Cells.Item(19, 6).Value = 558
$ws.Cells.Item(20, 6).Value = 446
$ws.Cells.Item(21, 6).Value = 462
$ws.Cells.Item(23, 6).Value = 11
$ws.Cells.Item(25, 6).Value = 855
$ws.Cells.Item(26, 6).Value = 122
$ws.Cells.Item(27, 6).Value = 37
$ws.Cells.Item(28, 6).Value = 8
$ws.Cells.Item(31, 6).Value = 15
$ws.Cells.Item(32, 6).Value = 13
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 6).Value = 53
$ws.Cells.Item(6, 6).Value = 86
$ws.Cells.Item(8, 6).Value = 1612
$ws.Cells.Item(15, 6).Value = 1651
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 6).Value = 229
$ws.Cells.Item(3, 6).Value = 30
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 229
$ws.Cells.Item(3, 6).Value = 895
$ws.Cells.Item(3, 7).Value = 70
$ws.Cells.Item(4, 6).Value = 13957
$ws.Cells.Item(4, 7).Value = 88
$ws.Cells.Item(5, 6).Value = 13790
$ws.Cells.Item(5, 7).Value = 90
$ws.Cells.Item(6, 6).Value = 1065
$ws.Cells.Item(6, 7).Value = 70
$ws.Cells.Item(8, 7).Value = "不可售"
$ws.Cells.Item(10, 7).Value = "不可售"
$ws.Cells.Item(13, 6).Value = 784
$ws.Cells.Item(14, 6).Value = 53
$ws.Cells.Item(17, 6).Value = 30
$ws.Cells.Item(18, 6).Value = 142
$ws.Cells.Item(19, 6).Value = 105
$ws.Cells.Item(20, 6).Value = 89
$ws.Cells.Item(21, 6).Value = 154
$ws.Cells.Item(24, 6).Value = 86
$ws.Cells.Item(26, 6).Value = 558
$ws.Cells.Item(27, 6).Value = 446
$ws.Cells.Item(28, 6).Value = 462
$ws.Cells.Item(30, 6).Value = 11
$ws.Cells.Item(32, 6).Value = 855
$ws.Cells.Item(34, 6).Value = 1612
$ws.Cells.Item(39, 6).Value = 122
$ws.Cells.Item(40, 6).Value = 37
$ws.Cells.Item(41, 6).Value = 8
$ws.Cells.Item(46, 6).Value = 15
$ws.Cells.Item(47, 6).Value = 13
$ws.Cells.Item(48, 6).Value = 1651
